$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B26").Value = "date-picker"
$ws.Range("A26").Value = "019_url_date_picker"

$ws.Range("A27").Value = "019_date_picker"
$ws.Range("B27").Value = "June 11 2001"

$ws.Range("A28").Value = "019_date_and_time_picker"
$ws.Range("B28").Value = "October 2 1879"

$ws.Range("F15").Select()
